$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully finish F2 (plain then rich) BEFORE touching F3/F4
$c = $ws.Range("F2")
$c.Value = "5511913592962;5511913592962"
$chars2 = $c.Characters(15, 13)
$chars2.Font.Name = "Calibri"
$chars2.Font.Size = 11
$chars2.Font.Color = 0

# Now F3 - same flattened text; will it dedupe to F2 (rich) or create new since F2 "used up" the plain version already?
$ws.Range("F3").Value = "5511913592962;5511913592962"

$ws.Range("F4").Value = "5511913592962;5511940280229"
